$wb = $excel.ActiveWorkbook

# 1. The "data" sheet no longer tracks the groomed_file column -- remove
#    column B (groomed_file / *_DT.nrrd values), leaving only segmentation_file.
$dataSheet = $wb.Worksheets.Item("data")
$dataSheet.Columns.Item(2).Delete()
$dataSheet.Range("B1").Select()

# 2. The project now opens on the "data" page (tool_state/view_state) instead
#    of the "groom" page -- update the saved studio state accordingly.
$studioSheet = $wb.Worksheets.Item("studio")
$studioSheet.Range("B3").Value = "data"
$studioSheet.Range("B4").Value = "Original"

# 3. Re-number the parameter-page sheets: recreate each one (content is
#    unchanged) so the workbook picks up fresh sheetIds, matching the new
#    project-file version stamp.
$names = @("groom", "optimize", "analysis", "studio")
for ($round = 0; $round -lt 2; $round++) {
    foreach ($nm in $names) {
        $src = $wb.Worksheets.Item($nm)
        $src.Copy($null, $src)
        $wb.Worksheets.Item($nm).Delete()
        $wb.Worksheets.Item("$nm (2)").Name = $nm
    }
}

$wb.Worksheets.Item("data").Activate()
